# Refresh the cryptocurrency price (column D) and 1h volume-change
# (column E) figures on the active sheet to the latest scraped values.
# Column D is stored as text in the source data (it can contain
# "thousand dot" groupings like "67.357.82" that are not valid
# numbers), so every D write is apostrophe-prefixed to force Excel
# to keep it as literal text instead of re-typing it as a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.357.82'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '''3.489.39'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''596.50'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").Value = '''180.02'
$ws.Range("E6").Value = '  +4.35%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '''0.607'
$ws.Range("E8").Value = '  +3.56%  '
$ws.Range("D9").Value = '''3.492.59'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +4.74%  '
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").Value = '''0.436'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '''4.093.15'
$ws.Range("D14").Value = '''32.20'
$ws.Range("E14").Value = '  +9.13%  '
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("D16").Value = '''67.367.74'
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '''3.500.43'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").Value = '''14.28'
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").Value = '''389.73'
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").Value = '''7.93'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '''73.97'
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("D24").Value = '''0.542'
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").Value = '''5.73'
$ws.Range("E26").Value = '  +0.84%  '
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("D28").Value = '''10.35'
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("E29").Value = '  -2.74%  '
$ws.Range("D30").Value = '''1.01'
$ws.Range("E30").Value = '  +1.95%  '
$ws.Range("D31").Value = '''6.18'
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").Value = '''2.06'
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("D34").Value = '''23.55'
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").Value = '''7.38'
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = '''163.25'
$ws.Range("E38").Value = '  +0.63%  '
$ws.Range("D39").Value = '''0.870'
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("E40").Value = '  +10.88%  '
$ws.Range("E41").Value = '  -0.85%  '
$ws.Range("E42").Value = '  -0.70%  '
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = '''2.846.49'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").Value = '''26.42'
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("D47").Value = '''0.0722'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("D48").Value = '''41.72'
$ws.Range("E48").Value = '  -2.26%  '
$ws.Range("D49").Value = '''0.0300'
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").Value = '''333.90'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("E51").Value = '  -1.35%  '
